$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from 45205
# (2023-10-06) to 45206 (2023-10-07) for every data row (rows 2 through 216).
for ($r = 2; $r -le 216; $r++) {
    $ws.Cells.Item($r, 3).Value = 45206
}
